$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Periodo Mora" rows (16-20) are reordered from ascending (2405..2409)
# to descending (2409..2405), while each period keeps its own "Valor Mora".
# Column E is formatted as Text, so values are assigned as strings.

$ws.Range("E16").Value = "2409"
$ws.Range("F16").Value = 70180

$ws.Range("E17").Value = "2408"
$ws.Range("F17").Value = 70180

$ws.Range("E18").Value = "2407"
$ws.Range("F18").Value = 70180

$ws.Range("E19").Value = "2406"
$ws.Range("F19").Value = 70180

$ws.Range("E20").Value = "2405"
$ws.Range("F20").Value = 56144
